$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- Rows 6-8: mark status as "Hecho" and register day-2 consumption ---
$ws.Range("F6").Value = "Hecho"
$ws.Range("F7").Value = "Hecho"
$ws.Range("F8").Value = "Hecho"

$ws.Range("T6").Value = 1
$ws.Range("T7").Value = 1
$ws.Range("T8").Value = 2

# --- Row 9: new task "Realizar descripción de CU 01 y CU 02" ---
$ws.Range("D9").Value = "Realizar descripción de CU 01 y CU 02"
$ws.Range("E9").Value = "Mario"
$ws.Range("F9").Value = "Por iniciar"
$ws.Rows.Item(9).RowHeight = 30

# --- Row 10: new task "Realizar descripción de CU 03 y CU 04" ---
$ws.Range("B10").ClearContents()
$ws.Range("D10").Value = "Realizar descripción de CU 03 y CU 04"
$ws.Range("E10").Value = "Victor"
$ws.Range("F10").Value = "Por iniciar"
$ws.Rows.Item(10).RowHeight = 30

# --- Row 11: new task "Realizar mockup de ventana principal de profesor" ---
$ws.Range("D11").Value = "Realizar mockup de ventana principal de profesor"
$ws.Range("E11").Value = "Mario"
$ws.Range("F11").Value = "Por iniciar"
$ws.Rows.Item(11).RowHeight = 45

# --- Row 12: new task "Realizar mockup de ventana principal de director" ---
$ws.Range("D12").Value = "Realizar mockup de ventana principal de director"
$ws.Range("E12").Value = "Victor"
$ws.Range("F12").Value = "Por iniciar"
$ws.Rows.Item(12).RowHeight = 45

# --- Update the selected / active cell ---
$ws.Range("T8").Select()
